# Add a "rating" column to shopee_products.xlsx
#
# Before:  A=url  B=name  C=original_price  D=current_price
#          E=description  F=image_url  G=shopid  H=itemid
#
# After:   A=url  B=name  C=original_price  D=current_price
#          E=description  F=rating  G=image_url  H=shopid  I=itemid
#
# i.e. insert a blank column before the current column F (image_url),
# add a "rating" header, and populate it with a per-row rating value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift image_url / shopid / itemid one column to the right. Excel carries
# the column-F header cell's style onto the freshly inserted blank F1, so
# it already comes out bold/bordered like the rest of row 1.
$ws.Columns.Item(6).Insert()

# Plain text header - not numeric-looking, so Value assigns it as a
# shared string with no extra formatting needed.
$ws.Cells.Item(1, 6).Value = "rating"

# Per-row rating values (rows 2..41), in sheet order.
$ratings = @(
    "4.7", "4.7", "4.8", "4.8", "4.8", "4.8", "4.9", "4.9", "4.8", "4.8",
    "4.9", "4.7", "4.9", "4.9", "4.9", "4.9", "4.7", "4.8", "4.9", "4.9",
    "4.8", "4.7", "4.8", "4.8", "4.9", "4.8", "4.9", "4.9", "4.9", "4.9",
    "4.9", "4.8", "4.9", "4.3", "4.9", "5",   "4.9", "4.9", "4.8", "4.8"
)

# These values look numeric ("4.7", "5", ...), so a plain Value assignment
# would be auto-coerced to a number (like typing 4.7 into Excel). The
# source workbook instead stores them as plain shared-string text (no
# special number format), matching how shopid/itemid are stored. Using a
# literal-string formula ("="4.7"") then flattening it to a static value
# via Copy/PasteSpecial(xlPasteValues) reproduces that: the destination
# cell ends up as plain text with no NumberFormat/style side effects.
$scratch = $ws.Cells.Item(1048576, 1)
for ($i = 0; $i -lt $ratings.Length; $i++) {
    $row = $i + 2
    $scratch.Formula = '="' + $ratings[$i] + '"'
    $scratch.Copy()
    $ws.Cells.Item($row, 6).PasteSpecial(-4163)
}
$scratch.Clear()

# Matches the author's last on-screen selection in the saved file.
$ws.Range("P16").Select() | Out-Null
